$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cd34"
$ws.Cells.Item(2, 3).Value = "Sele"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 146.901596
$ws.Cells.Item(2, 8).Value = 440.704788
$ws.Cells.Item(2, 9).Value = 0.5061978858527532
$ws.Cells.Item(2, 10).Value = 0.5061978858527532
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.699506666666667
$ws.Cells.Item(2, 14).Value = 14.09852
$ws.Cells.Item(2, 15).Value = 0.9660495246229048
$ws.Cells.Item(2, 16).Value = 0.9660495246229047
$ws.Cells.Item(2, 17).Value = 690.3650297459735
$ws.Cells.Item(2, 18).Value = 6213.285267713761
$ws.Cells.Item(2, 19).Value = 0.4890122269931717
$ws.Cells.Item(2, 20).Value = 0.4890122269931716

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cd34"
$ws.Cells.Item(3, 3).Value = "Sele"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 146.901596
$ws.Cells.Item(3, 8).Value = 440.704788
$ws.Cells.Item(3, 9).Value = 0.5061978858527532
$ws.Cells.Item(3, 10).Value = 0.5061978858527532
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.1651576666666667
$ws.Cells.Item(3, 14).Value = 0.495473
$ws.Cells.Item(3, 15).Value = 0.03395047537709522
$ws.Cells.Item(3, 16).Value = 0.03395047537709522
$ws.Cells.Item(3, 17).Value = 24.26192482496934
$ws.Cells.Item(3, 18).Value = 218.357323424724
$ws.Cells.Item(3, 19).Value = 0.01718565885958155
$ws.Cells.Item(3, 20).Value = 0.01718565885958155

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Cd34"
$ws.Cells.Item(4, 3).Value = "Sele"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 138.990916
$ws.Cells.Item(4, 8).Value = 416.972748
$ws.Cells.Item(4, 9).Value = 0.4789390295796214
$ws.Cells.Item(4, 10).Value = 0.4789390295796214
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.699506666666667
$ws.Cells.Item(4, 14).Value = 14.09852
$ws.Cells.Item(4, 15).Value = 0.9660495246229048
$ws.Cells.Item(4, 16).Value = 0.9660495246229047
$ws.Cells.Item(4, 17).Value = 653.1887363481068
$ws.Cells.Item(4, 18).Value = 5878.698627132961
$ws.Cells.Item(4, 19).Value = 0.4626788218487486
$ws.Cells.Item(4, 20).Value = 0.4626788218487486

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cd34"
$ws.Cells.Item(5, 3).Value = "Sele"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 138.990916
$ws.Cells.Item(5, 8).Value = 416.972748
$ws.Cells.Item(5, 9).Value = 0.4789390295796214
$ws.Cells.Item(5, 10).Value = 0.4789390295796214
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.1651576666666667
$ws.Cells.Item(5, 14).Value = 0.495473
$ws.Cells.Item(5, 15).Value = 0.03395047537709522
$ws.Cells.Item(5, 16).Value = 0.03395047537709522
$ws.Cells.Item(5, 17).Value = 22.95541537442267
$ws.Cells.Item(5, 18).Value = 206.598738369804
$ws.Cells.Item(5, 19).Value = 0.01626020773087282
$ws.Cells.Item(5, 20).Value = 0.01626020773087282

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Cd34"
$ws.Cells.Item(6, 3).Value = "Sele"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.164219
$ws.Cells.Item(6, 8).Value = 0.492657
$ws.Cells.Item(6, 9).Value = 0.0005658707112811305
$ws.Cells.Item(6, 10).Value = 0.0005658707112811305
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.699506666666667
$ws.Cells.Item(6, 14).Value = 14.09852
$ws.Cells.Item(6, 15).Value = 0.9660495246229048
$ws.Cells.Item(6, 16).Value = 0.9660495246229047
$ws.Cells.Item(6, 17).Value = 0.7717482852933334
$ws.Cells.Item(6, 18).Value = 6.945734567640001
$ws.Cells.Item(6, 19).Value = 0.0005466591316311611
$ws.Cells.Item(6, 20).Value = 0.0005466591316311611

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Cd34"
$ws.Cells.Item(7, 3).Value = "Sele"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.164219
$ws.Cells.Item(7, 8).Value = 0.492657
$ws.Cells.Item(7, 9).Value = 0.0005658707112811305
$ws.Cells.Item(7, 10).Value = 0.0005658707112811305
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.1651576666666667
$ws.Cells.Item(7, 14).Value = 0.495473
$ws.Cells.Item(7, 15).Value = 0.03395047537709522
$ws.Cells.Item(7, 16).Value = 0.03395047537709522
$ws.Cells.Item(7, 17).Value = 0.02712202686233333
$ws.Cells.Item(7, 18).Value = 0.244098241761
$ws.Cells.Item(7, 19).Value = 0.00001921157964996938
$ws.Cells.Item(7, 20).Value = 0.00001921157964996938

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Cd34"
$ws.Cells.Item(8, 3).Value = "Sele"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.149135333333334
$ws.Cells.Item(8, 8).Value = 12.447406
$ws.Cells.Item(8, 9).Value = 0.01429721385634429
$ws.Cells.Item(8, 10).Value = 0.01429721385634429
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.699506666666667
$ws.Cells.Item(8, 14).Value = 14.09852
$ws.Cells.Item(8, 15).Value = 0.9660495246229048
$ws.Cells.Item(8, 16).Value = 0.9660495246229047
$ws.Cells.Item(8, 17).Value = 19.49888915990223
$ws.Cells.Item(8, 18).Value = 175.49000243912
$ws.Cells.Item(8, 19).Value = 0.01381181664935341
$ws.Cells.Item(8, 20).Value = 0.01381181664935341

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Cd34"
$ws.Cells.Item(9, 3).Value = "Sele"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.149135333333334
$ws.Cells.Item(9, 8).Value = 12.447406
$ws.Cells.Item(9, 9).Value = 0.01429721385634429
$ws.Cells.Item(9, 10).Value = 0.01429721385634429
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.1651576666666667
$ws.Cells.Item(9, 14).Value = 0.495473
$ws.Cells.Item(9, 15).Value = 0.03395047537709522
$ws.Cells.Item(9, 16).Value = 0.03395047537709522
$ws.Cells.Item(9, 17).Value = 0.6852615103375557
$ws.Cells.Item(9, 18).Value = 6.167353593038
$ws.Cells.Item(9, 19).Value = 0.0004853972069908815
$ws.Cells.Item(9, 20).Value = 0.0004853972069908815
